$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.875.43"
$ws.Range("E2").Value = "  +1.50%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.622.28"
$ws.Range("E3").Value = "  +1.07%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.993"
$ws.Range("E4").Value = "  -0.73%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.50"
$ws.Range("E5").Value = "  +0.41%  "

# Row 6
$ws.Range("E6").Value = "  -0.28%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.72%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.13"
$ws.Range("E8").Value = "  +8.72%  "

# Row 9
$ws.Range("E9").Value = "  +3.31%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0607"
$ws.Range("E10").Value = "  +0.86%  "

# Row 11
$ws.Range("E11").Value = "  +0.03%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.94"
$ws.Range("E12").Value = "  +0.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.619.85"
$ws.Range("E13").Value = "  +0.76%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.568"
$ws.Range("E14").Value = "  +6.30%  "

# Row 15
$ws.Range("E15").Value = "  +5.57%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.902.54"
$ws.Range("E16").Value = "  +1.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.84"
$ws.Range("E17").Value = "  +16.15%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.41"
$ws.Range("E18").Value = "  +1.57%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.03"
$ws.Range("E19").Value = "  +0.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0707"
$ws.Range("E20").Value = "  +2.55%  "

# Row 21
$ws.Range("E21").Value = "  -0.59%  "

# Row 22
$ws.Range("E22").Value = "  +2.79%  "

# Row 23
$ws.Range("E23").Value = "  +4.34%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  +1.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.68"
$ws.Range("E25").Value = "  +0.78%  "

# Row 26
$ws.Range("E26").Value = "  +2.50%  "

# Row 27
$ws.Range("E27").Value = "  +1.24%  "

# Row 28
$ws.Range("E28").Value = "  +3.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.993"
$ws.Range("E29").Value = "  -0.77%  "

# Row 30
$ws.Range("E30").Value = "  +3.30%  "

# Row 31
$ws.Range("E31").Value = "  +5.27%  "

# Row 32
$ws.Range("E32").Value = "  +3.90%  "

# Row 33
$ws.Range("E33").Value = "  +4.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.417.52"
$ws.Range("E34").Value = "  +0.35%  "

# Row 35
$ws.Range("E35").Value = "  +7.10%  "

# Row 36
$ws.Range("E36").Value = "  -0.16%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.87"
$ws.Range("E37").Value = "  +1.45%  "

# Row 38
$ws.Range("E38").Value = "  -0.84%  "

# Row 39
$ws.Range("E39").Value = "  +2.60%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.555"
$ws.Range("E40").Value = "  +3.24%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0503"
$ws.Range("E41").Value = "  +3.53%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("E42").Value = "  -0.10%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.828"
$ws.Range("E43").Value = "  +4.03%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.22"
$ws.Range("E44").Value = "  +5.10%  "

# Row 45
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "53.41"
$ws.Range("E45").Value = "  +0.94%  "

# Row 46
$ws.Range("E46").Value = "  +18.56%  "

# Row 47
$ws.Range("E47").Value = "  -0.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.42"
$ws.Range("E48").Value = "  +2.80%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.762.41"
$ws.Range("E49").Value = "  +1.03%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.50"
$ws.Range("E50").Value = "  +2.17%  "

# Row 51
$ws.Range("E51").Value = "  +7.68%  "
